$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1464
$ws.Range("F3").Value = 1435
$ws.Range("F6").Value = 696
$ws.Range("F7").Value = 35
$ws.Range("F8").Value = 630
$ws.Range("F9").Value = 471
$ws.Range("F10").Value = 73
$ws.Range("F11").Value = 1380
$ws.Range("F12").Value = 32808
$ws.Range("F13").Value = 7036
$ws.Range("F15").Value = 362
$ws.Range("F16").Value = 571
$ws.Range("F17").Value = 439
$ws.Range("F19").Value = 103
$ws.Range("F21").Value = 45
$ws.Range("F22").Value = 447
$ws.Range("F23").Value = 103
$ws.Range("F24").Value = 793
$ws.Range("F25").Value = 6
$ws.Range("F26").Value = 317
$ws.Range("F27").Value = 386
$ws.Range("F28").Value = 439
$ws.Range("F30").Value = 203
$ws.Range("F31").Value = 49
$ws.Range("F32").Value = 738
$ws.Range("F33").Value = 290
$ws.Range("F35").Value = 737
$ws.Range("F36").Value = 111
$ws.Range("F38").Value = 790
$ws.Range("F39").Value = 288
$ws.Range("F41").Value = 23

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 1180
$ws.Range("F5").Value = 157
$ws.Range("F6").Value = 290
$ws.Range("F7").Value = 4323
$ws.Range("F9").Value = 232
$ws.Range("F11").Value = 5
$ws.Range("F17").Value = 155
$ws.Range("F19").Value = 4294

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1448
$ws.Range("F3").Value = 352

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1448
$ws.Range("F3").Value = 352
$ws.Range("F4").Value = 1180
$ws.Range("F5").Value = 1464
$ws.Range("F7").Value = 1435
$ws.Range("F9").Value = 696
$ws.Range("F10").Value = 35
$ws.Range("F11").Value = 630
$ws.Range("F13").Value = 1380
$ws.Range("F14").Value = 157
$ws.Range("F15").Value = 290
$ws.Range("F17").Value = 232
$ws.Range("F18").Value = 232
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 7036
$ws.Range("F23").Value = 362
$ws.Range("F25").Value = 571
$ws.Range("F26").Value = 439
$ws.Range("F28").Value = 103
$ws.Range("F31").Value = 447
$ws.Range("F32").Value = 103
$ws.Range("F33").Value = 793
$ws.Range("F34").Value = 317
$ws.Range("F35").Value = 386
$ws.Range("F36").Value = 439
$ws.Range("F38").Value = 203
$ws.Range("F39").Value = 49
$ws.Range("F40").Value = 738
$ws.Range("F42").Value = 290
$ws.Range("F44").Value = 111
$ws.Range("F45").Value = 790
$ws.Range("F46").Value = 288
$ws.Range("F49").Value = 23
